# Scheduled-runner update: refresh computed profit/price columns (H:N) across sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 1657.7446
$ws.Range("I15").Value = 1657.7446
$ws.Range("K15").Value = 4973.2338
$ws.Range("M15").Value = -4804.2338
$ws.Range("H33").Value = 3788946
$ws.Range("I33").Value = 1661
$ws.Range("J33").Value = 8264828
$ws.Range("K33").Value = 1661
$ws.Range("L33").Value = 8264828
$ws.Range("M33").Value = -1432
$ws.Range("N33").Value = -8265286
$ws.Range("H129").Value = 1004.2
$ws.Range("I129").Value = 547.0714
$ws.Range("J129").Value = 1210.6451
$ws.Range("K129").Value = 1641.2142
$ws.Range("L129").Value = 3631.9353
$ws.Range("M129").Value = 3358.7858
$ws.Range("N129").Value = -13631.9353
$ws.Range("H131").Value = 1583.3334
$ws.Range("I131").Value = 1583.3334
$ws.Range("J131").Value = 0
$ws.Range("K131").Value = 4750.0002
$ws.Range("L131").Value = 0
$ws.Range("M131").Value = 289.9997999999996
$ws.Range("N131").ClearContents()
$ws.Range("H132").Value = 1574.1538
$ws.Range("I132").Value = 1653.5652
$ws.Range("J132").Value = 965.3333
$ws.Range("K132").Value = 4960.6956
$ws.Range("L132").Value = 2895.9999
$ws.Range("M132").Value = -2430.6956
$ws.Range("N132").Value = -7955.9999
$ws.Range("H137").Value = 1852.0938
$ws.Range("I137").Value = 1899.9
$ws.Range("J137").Value = 1772.4166
$ws.Range("K137").Value = 5699.700000000001
$ws.Range("L137").Value = 5317.2498
$ws.Range("M137").Value = -3149.700000000001
$ws.Range("N137").Value = -10417.2498

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6530.4097
$ws.Range("I32").Value = 4636.037
$ws.Range("K32").Value = 4636.037
$ws.Range("M32").Value = -4349.037
$ws.Range("H45").Value = 17278.5
$ws.Range("I45").Value = 17278.5
$ws.Range("K45").Value = 17278.5
$ws.Range("M45").Value = -16901.5
$ws.Range("H74").Value = 1594.3334
$ws.Range("I74").Value = 1453.4117
$ws.Range("J74").Value = 1936.5714
$ws.Range("K74").Value = 1453.4117
$ws.Range("L74").Value = 1936.5714
$ws.Range("M74").Value = -579.4117000000001
$ws.Range("N74").Value = -3684.5714
$ws.Range("H77").Value = 1594.3334
$ws.Range("I77").Value = 1453.4117
$ws.Range("J77").Value = 1936.5714
$ws.Range("K77").Value = 7267.058500000001
$ws.Range("L77").Value = 9682.857
$ws.Range("M77").Value = -2899.058500000001
$ws.Range("N77").Value = -18418.857
$ws.Range("H97").Value = 1266.826
$ws.Range("I97").Value = 1070.625
$ws.Range("K97").Value = 1070.625
$ws.Range("M97").Value = -574.625
$ws.Range("H132").Value = 4040.9167
$ws.Range("I132").Value = 2034.6428
$ws.Range("J132").Value = 6849.7
$ws.Range("K132").Value = 6103.928400000001
$ws.Range("L132").Value = 20549.1
$ws.Range("M132").Value = -3573.928400000001
$ws.Range("N132").Value = -25609.1

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H26").Value = 19999.5
$ws.Range("I26").Value = 19999.5
$ws.Range("K26").Value = 19999.5
$ws.Range("M26").Value = -19707.5
$ws.Range("H28").Value = 30000
$ws.Range("J28").Value = 30000
$ws.Range("L28").Value = 30000
$ws.Range("N28").Value = -30588
$ws.Range("H80").Value = 285.13043
$ws.Range("I80").Value = 67.166664
$ws.Range("J80").Value = 362.05884
$ws.Range("K80").Value = 67.166664
$ws.Range("L80").Value = 362.05884
$ws.Range("M80").Value = 930.833336
$ws.Range("N80").Value = -2358.05884
$ws.Range("H83").Value = 285.13043
$ws.Range("I83").Value = 67.166664
$ws.Range("J83").Value = 362.05884
$ws.Range("K83").Value = 335.83332
$ws.Range("L83").Value = 1810.2942
$ws.Range("M83").Value = 4656.16668
$ws.Range("N83").Value = -11794.2942
$ws.Range("H86").Value = 1653.9524
$ws.Range("I86").Value = 1604.5714
$ws.Range("J86").Value = 1752.7142
$ws.Range("K86").Value = 1604.5714
$ws.Range("L86").Value = 1752.7142
$ws.Range("M86").Value = -481.5714
$ws.Range("N86").Value = -3998.7142
$ws.Range("H89").Value = 1653.9524
$ws.Range("I89").Value = 1604.5714
$ws.Range("J89").Value = 1752.7142
$ws.Range("K89").Value = 8022.857
$ws.Range("L89").Value = 8763.571
$ws.Range("M89").Value = -2406.857
$ws.Range("N89").Value = -19995.571
$ws.Range("H96").Value = 25214
$ws.Range("I96").Value = 10428
$ws.Range("J96").Value = 40000
$ws.Range("K96").Value = 10428
$ws.Range("L96").Value = 40000
$ws.Range("M96").Value = -7682
$ws.Range("N96").Value = -45492

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2524.2546
$ws.Range("I31").Value = 1727.2433
$ws.Range("J31").Value = 4162.5557
$ws.Range("K31").Value = 1727.2433
$ws.Range("L31").Value = 4162.5557
$ws.Range("M31").Value = -1432.2433
$ws.Range("N31").Value = -4752.5557
$ws.Range("H34").Value = 2524.2546
$ws.Range("I34").Value = 1727.2433
$ws.Range("J34").Value = 4162.5557
$ws.Range("K34").Value = 1727.2433
$ws.Range("L34").Value = 4162.5557
$ws.Range("M34").Value = -1525.2433
$ws.Range("N34").Value = -4566.5557

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 3740830.5
$ws.Range("I12").Value = 12500119
$ws.Range("K12").Value = 37500357
$ws.Range("M12").Value = -37500184
$ws.Range("H98").Value = 10000224
$ws.Range("I98").Value = 192.16667
$ws.Range("K98").Value = 576.50001
$ws.Range("M98").Value = 921.49999
$ws.Range("H105").Value = 6975
$ws.Range("J105").Value = 6975
$ws.Range("L105").Value = 20925
$ws.Range("N105").Value = -26167
$ws.Range("H133").Value = 39057.773
$ws.Range("I133").Value = 95526.45
$ws.Range("J133").Value = 8000
$ws.Range("K133").Value = 286579.35
$ws.Range("L133").Value = 24000
$ws.Range("M133").Value = -281519.35
$ws.Range("N133").Value = -34120

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 6030.9165
$ws.Range("I70").Value = 6536.0713
$ws.Range("J70").Value = 5323.7
$ws.Range("K70").Value = 6536.0713
$ws.Range("L70").Value = 5323.7
$ws.Range("M70").Value = -6266.0713
$ws.Range("N70").Value = -5863.7
$ws.Range("H73").Value = 6030.9165
$ws.Range("I73").Value = 6536.0713
$ws.Range("J73").Value = 5323.7
$ws.Range("K73").Value = 6536.0713
$ws.Range("L73").Value = 5323.7
$ws.Range("M73").Value = -5600.0713
$ws.Range("N73").Value = -7195.7
$ws.Range("H126").Value = 4803.048
$ws.Range("J126").Value = 2655.3547
$ws.Range("L126").Value = 7966.0641
$ws.Range("N126").Value = -12906.0641

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1156.9231
$ws.Range("J46").Value = 1400
$ws.Range("L46").Value = 1400
$ws.Range("N46").Value = -1776
$ws.Range("H122").Value = 4792963
$ws.Range("J122").Value = 1003340.5
$ws.Range("L122").Value = 3010021.5
$ws.Range("N122").Value = -3014921.5
$ws.Range("H132").Value = 15880678
$ws.Range("I132").Value = 30314214
$ws.Range("K132").Value = 90942642
$ws.Range("M132").Value = -90940112

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 2373.442
$ws.Range("J136").Value = 2031.2727
$ws.Range("L136").Value = 6093.8181
$ws.Range("N136").Value = -11193.8181
